$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (refresh from GitHub Actions job).
# Row 30/31 additionally swap coin identity (BitcoinCash <-> ImmutableX reordered).

$ws.Range("D2").Value = "'24.429.70"
$ws.Range("E2").Value = "'  -1.63%  "

$ws.Range("D3").Value = "'1.652.59"
$ws.Range("E3").Value = "'  -3.37%  "

$ws.Range("E4").Value = "'  -0.14%  "

$ws.Range("D5").Value = "'312.47"
$ws.Range("E5").Value = "'  +0.38%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.30%  "

$ws.Range("E7").Value = "'  -3.01%  "

$ws.Range("D8").Value = "'46.97"
$ws.Range("E8").Value = "'  -5.62%  "

$ws.Range("D9").Value = "'0.3256"
$ws.Range("E9").Value = "'  -5.70%  "

$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "'  -7.05%  "

$ws.Range("D11").Value = "'0.07041"
$ws.Range("E11").Value = "'  -6.87%  "

$ws.Range("E12").Value = "'  -0.17%  "

$ws.Range("D13").Value = "'5.934"
$ws.Range("E13").Value = "'  -6.14%  "

$ws.Range("D14").Value = "'19.42"
$ws.Range("E14").Value = "'  -8.04%  "

$ws.Range("D15").Value = "'6.589"
$ws.Range("E15").Value = "'  -6.78%  "

$ws.Range("D16").Value = "'1.652.92"
$ws.Range("E16").Value = "'  -3.59%  "

$ws.Range("D17").Value = "'0.00001045"
$ws.Range("E17").Value = "'  -8.10%  "

$ws.Range("D18").Value = "'0.06608"
$ws.Range("E18").Value = "'  -1.76%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.16%  "

$ws.Range("D20").Value = "'78.03"
$ws.Range("E20").Value = "'  -7.91%  "

$ws.Range("D21").Value = "'5.924"
$ws.Range("E21").Value = "'  -7.44%  "

$ws.Range("D22").Value = "'15.65"
$ws.Range("E22").Value = "'  -9.90%  "

$ws.Range("D23").Value = "'12.41"
$ws.Range("E23").Value = "'  -6.13%  "

$ws.Range("D24").Value = "'24.397.23"
$ws.Range("E24").Value = "'  -1.87%  "

$ws.Range("D25").Value = "'2.464"
$ws.Range("E25").Value = "'  +0.41%  "

$ws.Range("D26").Value = "'2.336"
$ws.Range("E26").Value = "'  -16.55%  "

$ws.Range("D27").Value = "'148.24"
$ws.Range("E27").Value = "'  -2.58%  "

$ws.Range("D28").Value = "'18.59"
$ws.Range("E28").Value = "'  -9.14%  "

$ws.Range("D29").Value = "'1.839.18"
$ws.Range("E29").Value = "'  -3.39%  "

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = "'124.11"
$ws.Range("E30").Value = "'  -6.63%  "

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.180"
$ws.Range("E31").Value = "'  -5.11%  "

$ws.Range("D32").Value = "'4.093"
$ws.Range("E32").Value = "'  -3.62%  "

$ws.Range("D33").Value = "'5.673"
$ws.Range("E33").Value = "'  -18.16%  "

$ws.Range("D34").Value = "'0.08435"
$ws.Range("E34").Value = "'  -4.23%  "

$ws.Range("D35").Value = "'1.663"
$ws.Range("E35").Value = "'  -9.93%  "

$ws.Range("D36").Value = "'12.30"
$ws.Range("E36").Value = "'  -11.21%  "

$ws.Range("D37").Value = "'5.203"
$ws.Range("E37").Value = "'  -7.50%  "

$ws.Range("D38").Value = "'0.06034"
$ws.Range("E38").Value = "'  -10.04%  "

$ws.Range("D39").Value = "'0.02215"
$ws.Range("E39").Value = "'  -8.31%  "

$ws.Range("D40").Value = "'0.2067"
$ws.Range("E40").Value = "'  -7.92%  "

$ws.Range("D41").Value = "'8.167"
$ws.Range("E41").Value = "'  -11.82%  "

$ws.Range("D42").Value = "'1.195"
$ws.Range("E42").Value = "'  -6.55%  "

$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "'  +0.34%  "

$ws.Range("D44").Value = "'0.5903"
$ws.Range("E44").Value = "'  -8.86%  "

$ws.Range("D45").Value = "'3.783"
$ws.Range("E45").Value = "'  -1.09%  "

$ws.Range("D46").Value = "'12.62"
$ws.Range("E46").Value = "'  -10.18%  "

$ws.Range("D47").Value = "'0.5622"
$ws.Range("E47").Value = "'  -9.19%  "

$ws.Range("D48").Value = "'122.39"
$ws.Range("E48").Value = "'  -5.98%  "

$ws.Range("D49").Value = "'1.942"
$ws.Range("E49").Value = "'  -9.47%  "

$ws.Range("D50").Value = "'0.06896"
$ws.Range("E50").Value = "'  -5.92%  "

$ws.Range("D51").Value = "'74.74"
$ws.Range("E51").Value = "'  -6.79%  "
